$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.332.27'
$ws.Range("E2").Value = '  +0.08%  '

$ws.Range("D3").Value = '1.789.44'
$ws.Range("E3").Value = '  +1.65%  '

$ws.Range("E4").Value = '  +1.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '337.33'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.37%  '

$ws.Range("E6").Value = '  +0.86%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3803'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.26%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3453'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.75%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.56'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.47%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.199'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.24%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07506'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.07%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.003'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.00%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.99'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.16%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.471'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.64%  '

$ws.Range("D15").Value = '1.790.40'
$ws.Range("E15").Value = '  +1.74%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.081'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.41%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001104'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.29%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06681'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.67%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '84.91'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.62%  '

$ws.Range("E20").Value = '  +0.91%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.567'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.74%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.40'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.92%  '

$ws.Range("D23").Value = '27.334.71'
$ws.Range("E23").Value = '  +0.39%  '

$ws.Range("E24").Value = '  -3.45%  '

$ws.Range("E25").Value = '  -0.64%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.511'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.14%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.570'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.00%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.62'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +9.09%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '153.37'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.56%  '

$ws.Range("D30").Value = '1.992.27'
$ws.Range("E30").Value = '  +1.93%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '134.01'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.58%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.058'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.60%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.086'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.35%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08701'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.59%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '13.20'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.51%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.648'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.72%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.491'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.19%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6924'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.99%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06412'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.05%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.873'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.06%  '

$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2203'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.07%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.02346'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.29%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.264'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.83%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.61'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.96%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6469'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.75%  '

$ws.Range("E46").Value = '  +0.70%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.863'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.14%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.138'
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '129.73'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.63%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07196'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.92%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.64'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.75%  '
